# Auto-generated script: apply 2022-11-02 data update to violent-crime-full-year.xlsx
# For each affected worksheet, update the specific cells per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 6161
$ws.Range('I3').Value = 6416
$ws.Range('F4').Value = 1869
$ws.Range('G4').Value = 1446
$ws.Range('I4').Value = 1474
$ws.Range('I5').Value = 593
$ws.Range('G6').Value = 7856
$ws.Range('I6').Value = 7285
$ws.Range('F7').Value = 24058
$ws.Range('G7').Value = 24671
$ws.Range('I7').Value = 21929

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I7').Value = 692
$ws.Range('I8').Value = 1316
$ws.Range('I9').Value = 111
$ws.Range('I11').Value = 329
$ws.Range('I15').Value = 252
$ws.Range('I16').Value = 64
$ws.Range('I19').Value = 605
$ws.Range('I23').Value = 215
$ws.Range('I25').Value = 118
$ws.Range('I27').Value = 193
$ws.Range('I29').Value = 1348
$ws.Range('I30').Value = 76
$ws.Range('I33').Value = 998
$ws.Range('I36').Value = 299
$ws.Range('I37').Value = 694
$ws.Range('I42').Value = 769
$ws.Range('I48').Value = 291
$ws.Range('I51').Value = 257
$ws.Range('I53').Value = 232
$ws.Range('I54').Value = 451
$ws.Range('I60').Value = 120
$ws.Range('F63').Value = 159
$ws.Range('G63').Value = 209
$ws.Range('I63').Value = 69
$ws.Range('I64').Value = 183
$ws.Range('I65').Value = 515
$ws.Range('I67').Value = 845
$ws.Range('I73').Value = 202
$ws.Range('I78').Value = 296
$ws.Range('I83').Value = 479
$ws.Range('I87').Value = 51
$ws.Range('I88').Value = 199
$ws.Range('I89').Value = 257
$ws.Range('I90').Value = 277
$ws.Range('I94').Value = 227
$ws.Range('I95').Value = 337
$ws.Range('I96').Value = 239
$ws.Range('G99').Value = 420
$ws.Range('I99').Value = 396
$ws.Range('F101').Value = 24058
$ws.Range('G101').Value = 24671
$ws.Range('I101').Value = 21929

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('I3').Value = 65
$ws.Range('I7').Value = 329

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I3').Value = 371
$ws.Range('I5').Value = 39
$ws.Range('I6').Value = 428
$ws.Range('I7').Value = 1316

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('I6').Value = 108
$ws.Range('I7').Value = 232

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('I2').Value = 227
$ws.Range('I4').Value = 37
$ws.Range('I6').Value = 182
$ws.Range('I7').Value = 692

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('I2').Value = 63
$ws.Range('I3').Value = 62
$ws.Range('I6').Value = 88
$ws.Range('I7').Value = 257

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('I6').Value = 89
$ws.Range('I7').Value = 239

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('I3').Value = 27
$ws.Range('I7').Value = 76

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I2').Value = 205
$ws.Range('I6').Value = 204
$ws.Range('I7').Value = 694

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I3').Value = 147
$ws.Range('G6').Value = 102
$ws.Range('I6').Value = 100
$ws.Range('G7').Value = 420
$ws.Range('I7').Value = 396

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I3').Value = 311
$ws.Range('I6').Value = 261
$ws.Range('I7').Value = 845

$ws = $wb.Worksheets.Item('New City')
$ws.Range('I2').Value = 171
$ws.Range('I7').Value = 515

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('I2').Value = 163
$ws.Range('I7').Value = 479

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('I2').Value = 116
$ws.Range('I7').Value = 337

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I2').Value = 226
$ws.Range('I3').Value = 371
$ws.Range('I4').Value = 45
$ws.Range('I6').Value = 316
$ws.Range('I7').Value = 998

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I2').Value = 97
$ws.Range('I6').Value = 216
$ws.Range('I7').Value = 451

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 397
$ws.Range('I6').Value = 374
$ws.Range('I7').Value = 1348

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I6').Value = 183
$ws.Range('I7').Value = 605

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('I3').Value = 58
$ws.Range('I6').Value = 150
$ws.Range('I7').Value = 291

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I2').Value = 192
$ws.Range('I6').Value = 261
$ws.Range('I7').Value = 769

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('I6').Value = 109
$ws.Range('I7').Value = 296

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('I2').Value = 60
$ws.Range('I7').Value = 215

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('I3').Value = 55
$ws.Range('I7').Value = 183

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('I4').Value = 11
$ws.Range('I7').Value = 299

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I6').Value = 129
$ws.Range('I7').Value = 227

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('I2').Value = 44
$ws.Range('I3').Value = 33
$ws.Range('I7').Value = 118

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('I3').Value = 58
$ws.Range('I7').Value = 252

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('I3').Value = 39
$ws.Range('I7').Value = 111

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('I2').Value = 66
$ws.Range('I7').Value = 202

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('I3').Value = 71
$ws.Range('I6').Value = 60
$ws.Range('I7').Value = 199

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('I6').Value = 75
$ws.Range('I7').Value = 193

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('I3').Value = 71
$ws.Range('I6').Value = 93
$ws.Range('I7').Value = 277

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('I3').Value = 69
$ws.Range('I6').Value = 104
$ws.Range('I7').Value = 257

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('I6').Value = 35
$ws.Range('I7').Value = 120

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('I4').Value = 6
$ws.Range('I7').Value = 51

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('I6').Value = 43
$ws.Range('I7').Value = 64
